# LabStarters/Lab03/Lab3Rubric-CIS195.xlsx
#
# Commit: "Added the rubric to the instructions"
#
# The two generic worksheet names are renamed to their real purpose
# (Sheet1 -> Rubric, Sheet2 -> Grade) now that the rubric workbook is
# linked from the lab instructions, and the "Rubric" sheet becomes the
# one the workbook opens on (active tab / selection) instead of "Grade".

$wb = $excel.ActiveWorkbook

$wsRubric = $wb.Worksheets.Item(1)
$wsGrade  = $wb.Worksheets.Item(2)

$wsRubric.Name = "Rubric"
$wsGrade.Name  = "Grade"

# Make "Rubric" the active/selected sheet (was "Grade" before the edit),
# with the last selection left on E14.
$wsRubric.Activate()
$wsRubric.Range("E14").Select()
